# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets
# to match freshly generated data (gh-pages output).

$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 2-4 -> column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 410
$ws1.Range("F3").Value = 2405
$ws1.Range("F4").Value = 114

# Sheet "全部类型": rows 2, 7, 8 -> column F (same events aggregated across types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 410
$ws4.Range("F7").Value = 2405
$ws4.Range("F8").Value = 114
